# ALASKA_2017.xlsx cleanup script
# - Renames header columns to clean machine-readable names
# - Normalizes "de"/"el"/"los" -> "De"/"El"/"Los" capitalization in a handful
#   of municipality / state names
# - Removes the trailing metadata/footer rows (77-81 and 476-480), which
#   shrinks the used range down to A1:D75

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header renames (row 1) ---
$ws.Range("A1").Value2 = "mx_state"
$ws.Range("B1").Value2 = "mx_municipality"
$ws.Range("C1").Value2 = "n_matriculas"
$ws.Range("D1").Value2 = "pct_matriculas"

# --- Capitalization fixes for particular place names ---
$ws.Range("A11").Value2 = "Ciudad De México"
$ws.Range("B20").Value2 = "Acapulco De Juárez"
$ws.Range("B21").Value2 = "Atoyac De Álvarez"
$ws.Range("B23").Value2 = "Cuautepec De Hinojosa"
$ws.Range("B24").Value2 = "Pachuca De Soto"
$ws.Range("B27").Value2 = "Tulancingo De Bravo"
$ws.Range("B33").Value2 = "San Juan De Los Lagos"
$ws.Range("B34").Value2 = "San Miguel El Alto"
$ws.Range("B36").Value2 = "Valle De Guadalupe"
$ws.Range("B46").Value2 = "Heroica Ciudad De Huajuapan De León"
$ws.Range("B47").Value2 = "Mariscala De Juárez"
$ws.Range("B48").Value2 = "Oaxaca De Juárez"
$ws.Range("B50").Value2 = "Tlacolula De Matamoros"
$ws.Range("B63").Value2 = "Amatlán De Los Reyes"
$ws.Range("B72").Value2 = "Tlaltenango De Sánchez Román"

# --- Minor recomputed percentage value ---
$ws.Range("D74").Value2 = 0.09195402298850576

# --- Remove trailing footer/metadata rows ---
# Delete from the bottom up so row numbers of the earlier block don't shift
$ws.Rows("476:480").Delete()
$ws.Rows("77:81").Delete()
